$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.378.77'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.81%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.651.77'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.33%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.01%  '

$ws.Range("E6").Value = '  +0.26%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3642'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.95%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '46.99'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.46%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3247'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.70%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.122'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.15%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07027'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -7.04%  '

$ws.Range("E12").Value = '  +0.09%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.935'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.79%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.37'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -8.32%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.588'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.61%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.651.09'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.44%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001042'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -8.18%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06605'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.72%  '

$ws.Range("E19").Value = '  +0.19%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '77.95'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -7.88%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.923'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -7.33%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.59'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -10.00%  '

$ws.Range("E23").Value = '  -6.67%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.363.70'
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.477'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.88%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.340'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -16.39%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '147.20'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.29%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.54'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -9.29%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.834.87'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.41%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.77'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -6.98%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.170'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.95%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.067'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.75%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.626'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -18.99%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08453'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.12%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.666'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.95%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.27'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -11.05%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.177'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -7.58%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.252'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.08%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06024'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -9.66%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02214'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -8.25%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2062'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -7.99%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.129'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -12.88%  '

$ws.Range("E43").Value = '  +0.20%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5879'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -8.90%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.768'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.40%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.59'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -10.39%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5604'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -9.28%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '122.39'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.83%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.943'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -9.17%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06879'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.90%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '74.55'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.81%  '
